$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a trailing space to the TODO note in F14 (text unchanged otherwise)
$ws.Range("F14").Value = "compare relevancies between courses| create a multi output regression model "

# Update hours worked on 2024-05-15 (row 14) from 4 to 8
$ws.Range("B14").Value = 8

# Move the active selection from F15 to B15
$ws.Range("B15").Select()
